# Update "想去人数" (F column) figures on both the "展览" and "全部类型" sheets.
# Both sheets share the same set of row updates (rows 3-38, skipping a few rows
# whose values did not change).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value (applies identically to both target sheets)
$updates = @{
    3  = 5600
    4  = 37
    5  = 678
    6  = 672
    7  = 27
    8  = 15
    10 = 228
    11 = 1584
    12 = 5426
    13 = 461
    14 = 289
    15 = 246
    16 = 55
    17 = 25
    19 = 4562
    20 = 232
    21 = 1211
    23 = 81
    24 = 214
    25 = 78
    26 = 202
    28 = 154
    29 = 83
    30 = 354
    31 = 38
    32 = 46
    34 = 19
    35 = 32
    37 = 44
    38 = 50
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
